# Update the "Robustness to window size" sheet with revised R2 values for
# the smaller window sizes (4 s, 2 s, 1 s, 0.5 s), and leave the cursor on
# the cell that was last edited, matching the author's final selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Robustness to window size")

$ws.Cells.Item(6, 2).Value = 0.8323
$ws.Cells.Item(7, 2).Value = 0.6789
$ws.Cells.Item(8, 2).Value = 0.2337
$ws.Cells.Item(9, 2).Value = 0.1106

$ws.Range("H11").Select()
